$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.465.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.87%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.958.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.79%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'198.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.77%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'596.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.27%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.04%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.203"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'2.959.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.77%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.445"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +11.13%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.504.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.71%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'28.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.74%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'76.473.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.99%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.961.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.16%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +8.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'378.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +4.80%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'72.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.05%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.093.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.66%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'4.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.47%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.77%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D31").Value = "'8.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +10.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.34%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'496.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.12%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'EthereumClassic"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'20.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.11%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Monero"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'164.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.12%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +14.71%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +19.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.53%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.04%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D43").Value = "'180.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.94%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'40.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.45%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.70%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.96%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.20%  "
$ws.Range("E51").Style = "Normal"
